# Delete the first data row (row 2, Trial_No 102) from Sheet1.
# Everything below shifts up by one row; the used range shrinks from
# A1:K115 to A1:K114.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows.Item(2).Delete()

# The AutoFilter range and the hidden _FilterDatabase defined name both
# still reference the old $K$115 extent after the row shift - refresh them
# to match the new used range.
$ws.AutoFilterMode = $false
$ws.Range("A1:K114").AutoFilter()
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$K`$114"

# Match the author's new selection recorded in the diff.
$ws.Range("N20").Select()
